$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 226, pushing existing rows 226-334 down to 227-335.
$ws.Rows("226").Insert()

# Populate the newly inserted row with the new weekly price record.
$ws.Range("A226").Value = 10
$ws.Range("B226").Value = "Vega Modelo de Temuco"
$ws.Range("C226").Value = "La Araucanía"
$ws.Range("D226").Value = 45134
$ws.Range("E226").Value = 9
$ws.Range("F226").Value = 100114007
$ws.Range("G226").Value = "Jengibre"
$ws.Range("H226").Value = "Sin especificar"
$ws.Range("I226").Value = "Primera"
$ws.Range("J226").Value = 35
$ws.Range("K226").Value = 24000
$ws.Range("L226").Value = 24000
$ws.Range("M226").Value = 24000
$ws.Range("N226").Value = "$/caja 13 kilos"
$ws.Range("O226").Value = "Perú"
$ws.Range("P226").Value = 1846
$ws.Range("Q226").Value = 13
$ws.Range("R226").Value = "Hortaliza"
